$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# ---------------------------------------------------------------------------
# 1) New row 76 / 77 values.
#    The shared-string table is built in the order cell values are written
#    (not in sheet row/column order), so the order below is deliberately
#    chosen to reproduce the exact shared-string indices from the target
#    workbook (Authoring75=269, Authoring76=270, OPQA-1195...=271,
#    "Verify saving..."=272, OPQA-1196...=273, "Verify draft title..."=274).
# ---------------------------------------------------------------------------
$ws.Range("A76").Value = "Authoring75"
$ws.Range("A77").Value = "Authoring76"
$ws.Range("B76").Value = "OPQA-1195|OPQA-1313|OPQA-1312|OPQA-1090|OPQA-1201"
$ws.Range("C76").Value = "Verify saving post as draft, accessing it for edit from profile,delete post from prfile"
$ws.Range("B77").Value = "OPQA-1196|OPQA-1200|OPQA-1199"
$ws.Range("C77").Value = "Verify draft title,access and edit draft post from post modal, delete post from post modal"

# D76 / D77 reuse the existing "Y" value already used throughout column D.
$ws.Range("D76").Value = $ws.Range("D75").Value2
$ws.Range("D77").Value = $ws.Range("D75").Value2

# ---------------------------------------------------------------------------
# 2) Formatting - copy the existing look-and-feel from cells that already
#    carry the desired format, so the new rows match the established
#    "Test Cases" table style.
# ---------------------------------------------------------------------------
$ws.Range("A75").Copy()
$ws.Range("A76").PasteSpecial(-4122)

$ws.Range("B64").Copy()
$ws.Range("B76").PasteSpecial(-4122)

$ws.Range("C42").Copy()
$ws.Range("C76").PasteSpecial(-4122)

$ws.Range("D72").Copy()
$ws.Range("D76").PasteSpecial(-4122)

$ws.Range("E42").Copy()
$ws.Range("E76").PasteSpecial(-4122)

$ws.Range("A75").Copy()
$ws.Range("A77").PasteSpecial(-4122)

$ws.Range("B61").Copy()
$ws.Range("B77").PasteSpecial(-4122)

$ws.Range("C42").Copy()
$ws.Range("C77").PasteSpecial(-4122)

$ws.Range("A75").Copy()
$ws.Range("D77").PasteSpecial(-4122)

$ws.Range("E42").Copy()
$ws.Range("E77").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3) D72:D75 lose their (redundant) fill flag once the sheet is re-saved.
# ---------------------------------------------------------------------------
$ws.Range("D72:D75").Interior.Pattern = -4142

# ---------------------------------------------------------------------------
# 4) Selection / scroll position.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 48
$ws.Range("D2:D75").Select()
